$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "MSME Number" column (column P) is being removed from the template.
# Deleting the whole column shifts row 1 (headers) and row 2 (placeholder
# formulas) one column to the left, and also shifts/reflows the column
# widths - exactly matching what Excel does when a user deletes a column.
$ws.Range("P1").EntireColumn.Delete()

# Row 1 header text fix: "PAN Number" -> "PAN"
$ws.Range("N1").Value = "PAN"

# Row 2 holds the merge-field placeholders. Before this edit row 2 was
# already one column "ahead" of row 1 (a pre-existing misalignment), so
# after the column delete above it must be restored to its correct,
# un-shifted values rather than left shifted like row 1.
$ws.Range("N2").Value = "{vendor:pan_no}"
$ws.Range("O2").Value = "{vendor:gst_no}"
$ws.Range("P2").Value = "{vendor:gst_taxpayer_type}"
$ws.Range("Q2").Value = "{vendor:gst_status}"
$ws.Range("R2").Value = "{vendor:appliances}"
$ws.Range("S2").Value = "{vendor:brands}"
$ws.Range("T2").Value = "{vendor:remarks}"
$ws.Range("U2").Value = "{vendor:primary_contact_name}"
$ws.Range("V2").Value = "{vendor:primary_contact_email}"
$ws.Range("W2").Value = "{vendor:primary_contact_phone_1}"
$ws.Range("X2").Value = "{vendor:primary_contact_phone_2}"
$ws.Range("Y2").Value = "{vendor:owner_name}"
$ws.Range("Z2").Value = "{vendor:owner_email}"
$ws.Range("AA2").Value = "{vendor:owner_phone_1}"
$ws.Range("AB2").Value = "{vendor:owner_phone_2}"
$ws.Range("AC2").Value = "{vendor:sf_create_date}"
$ws.Range("AD2").Value = "{vendor:non_working_days}"
$ws.Range("AE2").Value = "{vendor:sc_code}"
$ws.Range("AF2").Value = "{vendor:sf_rm_name}"
$ws.Range("AG2").Value = "{vendor:sf_rm_phone}"
$ws.Range("AH2").Value = "{vendor:sf_asm_name}"
$ws.Range("AI2").Value = "{vendor:sf_asm_phone}"
$ws.Range("AJ2").Value = "{vendor:active_status}"
$ws.Range("AK2").Value = "{vendor:on_off_status}"

# Restore the view/selection to show the corrected "Create Date" column.
$ws.Range("AC2").Select()
$excel.ActiveWindow.ScrollColumn = 23

Write-Output "done"
